# Refresh the hard-coded Katalon payment test-run results: the rows that
# feed the "Gateway Production" and "Upgrade Live Processors" test cases/
# suites get a fresh execution timestamp in column B ("Date"), and several
# of them now report "Fail" in column A ("Result") instead of "Pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ Result = <new Result or $null to leave unchanged>; Date = <new Date> }
$updates = [ordered]@{
    2  = @{ Date = "Tue Jul 11 12:27:51 EDT 2023" }
    3  = @{ Date = "Tue Jul 11 12:28:39 EDT 2023" }
    4  = @{ Date = "Tue Jul 11 12:29:28 EDT 2023" }
    5  = @{ Result = "Fail"; Date = "Tue Jul 11 12:30:11 EDT 2023" }
    6  = @{ Result = "Fail"; Date = "Tue Jul 11 12:30:53 EDT 2023" }
    7  = @{ Result = "Fail"; Date = "Tue Jul 11 12:31:34 EDT 2023" }
    8  = @{ Result = "Fail"; Date = "Tue Jul 11 12:32:16 EDT 2023" }
    9  = @{ Result = "Fail"; Date = "Tue Jul 11 12:33:04 EDT 2023" }
    10 = @{ Result = "Fail"; Date = "Tue Jul 11 12:33:53 EDT 2023" }
    11 = @{ Date = "Tue Jul 11 12:34:41 EDT 2023" }
    12 = @{ Date = "Tue Jul 11 12:35:25 EDT 2023" }
    13 = @{ Date = "Tue Jul 11 12:36:14 EDT 2023" }
    14 = @{ Result = "Fail"; Date = "Tue Jul 11 12:36:56 EDT 2023" }
}

foreach ($row in $updates.Keys) {
    $info = $updates[$row]
    if ($info.ContainsKey("Result")) {
        $ws.Cells.Item($row, 1).Value = $info.Result
    }
    $ws.Cells.Item($row, 2).Value = $info.Date
}

$wb.Save()
